# Natmi following Dr Hou advice
# Re-derive the Lama2-Rpsa ligand/receptor NATMI edge table so that the
# "ECs" (endothelial cells) cluster is included as both a sending and a
# target cluster alongside the existing "FAPs" and "sCs" clusters. This
# expands the previous 2 (senders) x 3 (targets) = 6 data rows into a full
# 3 x 3 = 9 data row grid (rows 2-10) and refreshes every computed NATMI
# statistic (detection rates, average/total expression, specificities,
# edge weights) to match the new 3-cluster denominator.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Lama2"
$ws.Cells.Item(2,3).Value = "Rpsa"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 2.452389333333333
$ws.Cells.Item(2,8).Value = 7.357168
$ws.Cells.Item(2,9).Value = 0.007993767302975028
$ws.Cells.Item(2,10).Value = 0.007993767302975028
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 52.80829433333333
$ws.Cells.Item(2,14).Value = 158.424883
$ws.Cells.Item(2,15).Value = 0.1724060238174878
$ws.Cells.Item(2,16).Value = 0.1724060238174878
$ws.Cells.Item(2,17).Value = 129.5064977345938
$ws.Cells.Item(2,18).Value = 1165.558479611344
$ws.Cells.Item(2,19).Value = 0.001378173636028168
$ws.Cells.Item(2,20).Value = 0.001378173636028168

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Lama2"
$ws.Cells.Item(3,3).Value = "Rpsa"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 2.452389333333333
$ws.Cells.Item(3,8).Value = 7.357168
$ws.Cells.Item(3,9).Value = 0.007993767302975028
$ws.Cells.Item(3,10).Value = 0.007993767302975028
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 81.07766966666667
$ws.Cells.Item(3,14).Value = 243.233009
$ws.Cells.Item(3,15).Value = 0.2646985445010758
$ws.Cells.Item(3,16).Value = 0.2646985445010758
$ws.Cells.Item(3,17).Value = 198.8340122620569
$ws.Cells.Item(3,18).Value = 1789.506110358512
$ws.Cells.Item(3,19).Value = 0.00211593857017778
$ws.Cells.Item(3,20).Value = 0.00211593857017778

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Lama2"
$ws.Cells.Item(4,3).Value = "Rpsa"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 2.452389333333333
$ws.Cells.Item(4,8).Value = 7.357168
$ws.Cells.Item(4,9).Value = 0.007993767302975028
$ws.Cells.Item(4,10).Value = 0.007993767302975028
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 172.4159456666667
$ws.Cells.Item(4,14).Value = 517.247837
$ws.Cells.Item(4,15).Value = 0.5628954316814363
$ws.Cells.Item(4,16).Value = 0.5628954316814364
$ws.Cells.Item(4,17).Value = 422.8310260495128
$ws.Cells.Item(4,18).Value = 3805.479234445616
$ws.Cells.Item(4,19).Value = 0.004499655096769079
$ws.Cells.Item(4,20).Value = 0.00449965509676908

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Lama2"
$ws.Cells.Item(5,3).Value = "Rpsa"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 243.5672963333334
$ws.Cells.Item(5,8).Value = 730.7018890000001
$ws.Cells.Item(5,9).Value = 0.7939278902575405
$ws.Cells.Item(5,10).Value = 0.7939278902575405
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 52.80829433333333
$ws.Cells.Item(5,14).Value = 158.424883
$ws.Cells.Item(5,15).Value = 0.1724060238174878
$ws.Cells.Item(5,16).Value = 0.1724060238174878
$ws.Cells.Item(5,17).Value = 12862.37347474489
$ws.Cells.Item(5,18).Value = 115761.361272704
$ws.Cells.Item(5,19).Value = 0.1368779507571094
$ws.Cells.Item(5,20).Value = 0.1368779507571094

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Lama2"
$ws.Cells.Item(6,3).Value = "Rpsa"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 243.5672963333334
$ws.Cells.Item(6,8).Value = 730.7018890000001
$ws.Cells.Item(6,9).Value = 0.7939278902575405
$ws.Cells.Item(6,10).Value = 0.7939278902575405
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 81.07766966666667
$ws.Cells.Item(6,14).Value = 243.233009
$ws.Cells.Item(6,15).Value = 0.2646985445010758
$ws.Cells.Item(6,16).Value = 0.2646985445010758
$ws.Cells.Item(6,17).Value = 19747.86879371711
$ws.Cells.Item(6,18).Value = 177730.819143454
$ws.Cells.Item(6,19).Value = 0.2101515569899808
$ws.Cells.Item(6,20).Value = 0.2101515569899808

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Lama2"
$ws.Cells.Item(7,3).Value = "Rpsa"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 243.5672963333334
$ws.Cells.Item(7,8).Value = 730.7018890000001
$ws.Cells.Item(7,9).Value = 0.7939278902575405
$ws.Cells.Item(7,10).Value = 0.7939278902575405
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 172.4159456666667
$ws.Cells.Item(7,14).Value = 517.247837
$ws.Cells.Item(7,15).Value = 0.5628954316814363
$ws.Cells.Item(7,16).Value = 0.5628954316814364
$ws.Cells.Item(7,17).Value = 41994.8857307849
$ws.Cells.Item(7,18).Value = 377953.9715770641
$ws.Cells.Item(7,19).Value = 0.4468983825104503
$ws.Cells.Item(7,20).Value = 0.4468983825104503

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Lama2"
$ws.Cells.Item(8,3).Value = "Rpsa"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 60.76799533333334
$ws.Cells.Item(8,8).Value = 182.303986
$ws.Cells.Item(8,9).Value = 0.1980783424394845
$ws.Cells.Item(8,10).Value = 0.1980783424394845
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 52.80829433333333
$ws.Cells.Item(8,14).Value = 158.424883
$ws.Cells.Item(8,15).Value = 0.1724060238174878
$ws.Cells.Item(8,16).Value = 0.1724060238174878
$ws.Cells.Item(8,17).Value = 3209.054183609293
$ws.Cells.Item(8,18).Value = 28881.48765248364
$ws.Cells.Item(8,19).Value = 0.03414989942435027
$ws.Cells.Item(8,20).Value = 0.03414989942435027

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Lama2"
$ws.Cells.Item(9,3).Value = "Rpsa"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 60.76799533333334
$ws.Cells.Item(9,8).Value = 182.303986
$ws.Cells.Item(9,9).Value = 0.1980783424394845
$ws.Cells.Item(9,10).Value = 0.1980783424394845
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 81.07766966666667
$ws.Cells.Item(9,14).Value = 243.233009
$ws.Cells.Item(9,15).Value = 0.2646985445010758
$ws.Cells.Item(9,16).Value = 0.2646985445010758
$ws.Cells.Item(9,17).Value = 4926.927451941542
$ws.Cells.Item(9,18).Value = 44342.34706747387
$ws.Cells.Item(9,19).Value = 0.05243104894091721
$ws.Cells.Item(9,20).Value = 0.05243104894091721

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Lama2"
$ws.Cells.Item(10,3).Value = "Rpsa"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 60.76799533333334
$ws.Cells.Item(10,8).Value = 182.303986
$ws.Cells.Item(10,9).Value = 0.1980783424394845
$ws.Cells.Item(10,10).Value = 0.1980783424394845
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 172.4159456666667
$ws.Cells.Item(10,14).Value = 517.247837
$ws.Cells.Item(10,15).Value = 0.5628954316814363
$ws.Cells.Item(10,16).Value = 0.5628954316814364
$ws.Cells.Item(10,17).Value = 10477.37138166425
$ws.Cells.Item(10,18).Value = 94296.34243497829
$ws.Cells.Item(10,19).Value = 0.111497394074217
$ws.Cells.Item(10,20).Value = 0.111497394074217

